$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: cell-level value edits on the original row numbering ---

# D3: -14.2 -> missing
$ws.Range("D3").ClearContents()

# E5: -5 -> missing
$ws.Range("E5").ClearContents()

# E8: missing -> -6.6
$ws.Range("E8").Value = -6.6

# E10: missing -> -6.1
$ws.Range("E10").Value = -6.1

# E12: -5.3 -> missing
$ws.Range("E12").ClearContents()

# E15: missing -> -8.4
$ws.Range("E15").Value = -8.4

# E18: -8.5 -> missing
$ws.Range("E18").ClearContents()

# E19: -6.5 -> missing
$ws.Range("E19").ClearContents()

# E25: missing -> -7.1
$ws.Range("E25").Value = -7.1

# --- Step 2: remove row 26 ("RM 232") entirely, shifting rows up ---
$ws.Rows(26).Delete()

# --- Step 3: remove the row now holding "SC 92" (originally row 28) entirely ---
$ws.Rows(27).Delete()

# --- Step 4: cell-level edits on the new row numbering (after the two deletions) ---

# Row 26 is now "SC 5": C26 missing -> 10.8
$ws.Range("C26").Value = 10.8

# Row 27 is now "SC 101": C27 10 -> missing
$ws.Range("C27").ClearContents()

# Row 29 is now "SC 119": E29 -6.8 -> missing
$ws.Range("E29").ClearContents()

# Row 33 is now "SC 232": C33 missing -> 10.4, D33 missing -> -14.1
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = -14.1
